# Updates crypto price/volume cells in the active sheet to match the
# latest scrape (GitHub Actions "Updated cryptos list" commit).
# D<row> holds the Price (plain text, e.g. "27.141.76" or "1.004"),
# E<row> holds the 1h Volume/% change (plain text with padding, e.g. "  -0.04%  ").
# Both columns are stored as TEXT in the workbook, never as numbers, so for any
# new Price value that Excel would otherwise auto-convert to a number (e.g.
# "1.004", "21.40"), the cell is temporarily forced to Text format, written,
# and then restored to the default "Normal" style so no formatting/style diff
# is introduced - only the text content changes, exactly like the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.141.76"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.900.51"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5234"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07290"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9037"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08202"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "1.862.33"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008661"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "27.183.48"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.129"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "2.118.83"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.452"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.316"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.821"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.903"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05046"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7941"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.224"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.363"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.631"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5734"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01989"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.132"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.618"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4895"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "38.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
